# Updated cryptos list with latest price / volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted figures (e.g. "331.70",
# "  +7.82%  ") where trailing zeros / exact spacing matter. Force the cells to
# text BEFORE writing so Excel does not silently coerce numeric-looking strings
# (e.g. "331.70" -> 331.7) to the Number type, then restore the default style.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.955.89"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.673.88"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "331.70"
$ws.Range("E5").Value = "  +7.82%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.3655"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "47.21"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "0.3221"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "1.144"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").Value = "0.07141"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "6.087"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "19.66"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "1.670.55"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "6.647"
$ws.Range("D17").Value = "0.00001049"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "0.06547"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "78.82"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").Value = "15.86"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "5.908"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "12.83"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").Value = "24.942.59"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "2.398"
$ws.Range("E26").Value = "  +4.16%  "
$ws.Range("D27").Value = "148.41"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "18.70"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "1.854.00"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "125.83"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Value = "1.184"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "4.085"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "5.791"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("D34").Value = "0.08478"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "1.661"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "5.158"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "0.06058"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02231"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").Value = "0.2088"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "8.242"
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "0.5965"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").Value = "13.72"
$ws.Range("E45").Value = "  +7.80%  "
$ws.Range("D46").Value = "3.848"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "0.5724"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "124.34"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "0.07005"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").Value = "  +3.89%  "

$textRange.Style = "Normal"
